# 6 hours by turn fix
# The teacher's weekly schedule shifted so that each class period starts
# 20 minutes earlier after the mid-morning break, lunch moved from 13:00 to
# 12:20, and three new afternoon time slots (17:30, 18:20) were appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full target grid (rows 2..17, columns A..F) reflecting the
# schedule after the fix.
$data = @(
    @("7:00",  "-", "-", "-", "-", "-"),
    @("7:50",  "-", "-", "-", "MCT-1A-Circuitos Elétricos", "MCT-1A-Circuitos Elétricos"),
    @("8:40",  "-", "ELT-2A-Acionamentos", "-", "-", "-"),
    @("9:30",  "Intervalo", "Intervalo", "Intervalo", "Intervalo", "Intervalo"),
    @("9:50",  "-", "ELT-2A-Acionamentos", "MEC-1A-Circuitos Elétricos", "ELT-1A-Circuitos Elétricos", "ELT-1A-Circuitos Elétricos"),
    @("10:40", "-", "MCT-2A-Acionamentos", "MEC-1A-Circuitos Elétricos", "-", "-"),
    @("11:30", "-", "-", "-", "-", "-"),
    @("12:20", "Almoço", "Almoço", "Almoço", "Almoço", "Almoço"),
    @("13:00", "-", "-", "-", "-", "-"),
    @("13:50", "-", "-", "-", "-", "-"),
    @("14:40", "-", "-", "-", "-", "-"),
    @("15:30", "Intervalo", "Intervalo", "Intervalo", "Intervalo", "Intervalo"),
    @("15:50", "-", "-", "-", "-", "-"),
    @("16:40", "-", "-", "-", "-", "-"),
    @("17:30", "-", "-", "-", "-", "-")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $row++
}

# Final row (18:20) only has a time label; the remaining cells stay blank.
$ws.Cells.Item(17, 1).Value = "18:20"
$ws.Cells.Item(17, 2).Value = ""
$ws.Cells.Item(17, 3).Value = ""
$ws.Cells.Item(17, 4).Value = ""
$ws.Cells.Item(17, 5).Value = ""
$ws.Cells.Item(17, 6).Value = ""
